$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from SCD0199 to SCD0011
$ws.Name = "SCD0011"

# Update the TC_ID cell (B2) from "DGS-214" to "SCD0011-030"
$ws.Range("B2").Value = "SCD0011-030"

# Column B needs to widen to fit the new, longer TC_ID text
$ws.Columns("B").ColumnWidth = 11.67

# Select cell B3 as the active cell (matches the recorded selection state)
$ws.Range("B3").Select()
